$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.671.05"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "2.278.64"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +9.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.49%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.23"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.80"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.32%  "

$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").Value = "2.620.00"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.866"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("D17").Value = "2.280.06"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").Value = "43.458.52"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("E19").Value = "  -1.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +11.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("E22").Value = "  -4.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.46"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.92"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.54%  "

$ws.Range("E30").Value = "  +1.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.15"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.60%  "

$ws.Range("E32").Value = "  -1.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0911"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  -5.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0349"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.67%  "

$ws.Range("E38").Value = "  -5.68%  "

$ws.Range("E39").Value = "  +4.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +17.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.40"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +13.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.41"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.49%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("E44").Value = "  +16.37%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.66"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("E48").Value = "  -2.07%  "

$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.45"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.18%  "

$ws.Range("E51").Value = "  +2.87%  "
